# Commit "Add files via upload": the document title run changes from
#   "Manual de Mantenimiento — WAYKISAFE"
# to
#   "Manual de Mantenimiento — WaykiSafeWAYKISAFE"
# i.e. "WaykiSafe" is inserted immediately before the existing "WAYKISAFE".
$d = $word.ActiveDocument

$emDash = [char]0x2014
$titleRange = $d.Paragraphs.Item(1).Range

$titleRange.Find.Execute(
    "Manual de Mantenimiento $emDash WAYKISAFE",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Manual de Mantenimiento $emDash WaykiSafeWAYKISAFE", 2
)
